# Insert a new data row before the current row 182, shifting existing
# rows 182-192 down to 183-193 (so D193 keeps the value that used to be
# in D192, etc.), then populate the newly inserted row 182 with the new
# record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(182).Insert()

$ws.Cells.Item(182, 1).Value = 5
$ws.Cells.Item(182, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(182, 3).Value = "Maule"
$ws.Cells.Item(182, 4).Value = 44516
$ws.Cells.Item(182, 5).Value = 7
$ws.Cells.Item(182, 6).Value = 100112003
$ws.Cells.Item(182, 7).Value = "Ajo"
$ws.Cells.Item(182, 8).Value = "Chino"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 200
$ws.Cells.Item(182, 11).Value = 20000
$ws.Cells.Item(182, 12).Value = 20000
$ws.Cells.Item(182, 13).Value = 20000
$ws.Cells.Item(182, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(182, 15).Value = "China"
$ws.Cells.Item(182, 16).Value = 2000
$ws.Cells.Item(182, 17).Value = 10
$ws.Cells.Item(182, 18).Value = "Hortaliza"
